$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.432.43"
$ws.Range("E2").Value = "  +3.77%  "

$ws.Range("D3").Value = "2.065.41"
$ws.Range("E3").Value = "  +5.99%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.28%  "

$ws.Range("E6").Value = "  +4.48%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.67"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.48%  "

$ws.Range("E9").Value = "  +5.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0759"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.101"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.70%  "

$ws.Range("D13").Value = "2.370.36"
$ws.Range("E13").Value = "  +6.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.774"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.32%  "

$ws.Range("D18").Value = "2.069.42"
$ws.Range("E18").Value = "  +5.84%  "

$ws.Range("D19").Value = "37.591.78"
$ws.Range("E19").Value = "  +4.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +25.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.76%  "

$ws.Range("D22").Value = "0.0₃0809"
$ws.Range("E22").Value = "  +3.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "225.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.73%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("E25").Value = "  +7.01%  "

$ws.Range("E26").Value = "  +3.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.87%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.89%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.128"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.32%  "

$ws.Range("E32").Value = "  +2.61%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.51%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +16.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0624"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.50%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("E39").Value = "  +8.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +17.93%  "

$ws.Range("E41").Value = "  -0.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +32.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0955"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.33%  "

$ws.Range("D44").Value = "1.464.47"
$ws.Range("E44").Value = "  +5.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.67%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.83%  "

$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.37%  "

$ws.Range("E49").Value = "  +5.92%  "

$ws.Range("E50").Value = "  +9.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.74%  "
